# Update the "Panel B (E-mini Futures)" rows (Avg Daily Volume, Diff_Vol, # Obs)
# in the Post-ZLB daily volume table with revised statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: Avg Daily Volume (Emini)
$ws.Range("D26").Value  = 1156019.258064516
$ws.Range("E26").Value  = 470721.7044122966
$ws.Range("F26").Value  = 943247.5
$ws.Range("G26").Value  = 1214222.5
$ws.Range("H26").Value  = 1408225
$ws.Range("I26").Value  = 62
$ws.Range("J26").Value  = 1395887.725806452
$ws.Range("K26").Value  = 421848.9578408938
$ws.Range("L26").Value  = 1143188.5
$ws.Range("M26").Value  = 1395650
$ws.Range("N26").Value  = 1659552.5
$ws.Range("O26").Value  = 62
$ws.Range("P26").Value  = 1457891.338709677
$ws.Range("Q26").Value  = 371585.8271050477
$ws.Range("R26").Value  = 1216031.25
$ws.Range("S26").Value  = 1410431
$ws.Range("T26").Value  = 1771407
$ws.Range("U26").Value  = 62
$ws.Range("V26").Value  = 1470363.774193548
$ws.Range("W26").Value  = 661735.5641658152
$ws.Range("X26").Value  = 1229644.5
$ws.Range("Y26").Value  = 1442094.5
$ws.Range("Z26").Value  = 1752537.5
$ws.Range("AA26").Value = 62
$ws.Range("AB26").Value = 1404754.370967742
$ws.Range("AC26").Value = 526276.9506933801
$ws.Range("AD26").Value = 1101111.5
$ws.Range("AE26").Value = 1377275
$ws.Range("AF26").Value = 1757069.5
$ws.Range("AG26").Value = 62

# Row 27: Diff_Vol (Ann - Day) (Emini) -- only the Mean column per period changes
$ws.Range("D27").Value  = 301872.0806451613
$ws.Range("J27").Value  = 62003.6129032258
$ws.Range("V27").Value  = -12472.43548387097
$ws.Range("AB27").Value = 53136.96774193548

# Row 28: # Obs (Emini) -- only the Mean column per period changes
$ws.Range("D28").Value  = 62
$ws.Range("J28").Value  = 62
$ws.Range("P28").Value  = 62
$ws.Range("V28").Value  = 62
$ws.Range("AB28").Value = 62
